# Update "想去人数" (number of people interested) figures on the
# "展览" (Exhibitions) sheet and the corresponding rows on the
# "全部类型" (All Types) aggregate sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 5     # 合肥·11.16合肥耽美同人only : 3 -> 5
$wsExpo.Range("F5").Value = 3806  # 合肥·第九届环形宇宙动漫游戏嘉年华 : 3799 -> 3806
$wsExpo.Range("F8").Value = 221   # 合肥·心动恋章·冬日序国乙&代号鸢同人only : 218 -> 221
$wsExpo.Range("F9").Value = 8     # 合肥·皖萌次元青年文化节 : 5 -> 8

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 5      # 合肥·11.16合肥耽美同人only : 3 -> 5
$wsAll.Range("F9").Value = 3806   # 合肥·第九届环形宇宙动漫游戏嘉年华 : 3799 -> 3806
$wsAll.Range("F13").Value = 221   # 合肥·心动恋章·冬日序国乙&代号鸢同人only : 218 -> 221
$wsAll.Range("F14").Value = 8     # 合肥·皖萌次元青年文化节 : 5 -> 8

$wb.Save()
